$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

$ws.Range("A2").Value = "T01_EstabLogin"
$ws.Range("A3").Value = "T02_Normal"
$ws.Range("A4").Value = "T03_EstabEntry"
$ws.Range("A5").Value = "T04_Tourism"

$ws.Range("A3").Select()
